$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the single remaining data row (row 2) with the new values
$ws.Range("A2").Value = "MSKU"
$ws.Range("B2").Value = 505392
$ws.Range("D2").Value = "A36"
$ws.Range("E2").Value = "SAN ANTONIO TERMINAL INTERNACIONAL S.A."
$ws.Range("F2").Value = "SI"
$ws.Range("G2").Value = "NO"
$ws.Range("H2").Value = "SI"
$ws.Range("I2").Value = "04/06/2024  22:00 "
$ws.Range("J2").Value = "28/05/2024 00:23"
$ws.Range("K2").Value = "MAERSK LONDRINA"
$ws.Range("L2").Value = "27/05/2024 13:00"

# Remove rows 3 through 7, which no longer exist in the factibilidad model
$ws.Range("A3:L7").Delete()
